$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The target cells hold numeric-looking values that are stored as text
# (t="inlineStr") in the workbook, so we must write them back as text
# rather than letting Excel auto-convert them to numbers. Prefixing with
# a single quote forces Excel to keep the entry as literal text, exactly
# as it would if a user typed '4.78 into the cell.
$ws.Range("B2").Value = "'4.78"
$ws.Range("B3").Value = "'3.49"
$ws.Range("B4").Value = "'3.52"
$ws.Range("B5").Value = "'3.89"
$ws.Range("B6").Value = "'4.22"
